$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 151, shifting the existing
# rows 151-178 down to 153-180.
$ws.Rows.Item(151).Insert()
$ws.Rows.Item(151).Insert()

# New row 151
$ws.Cells.Item(151, 1).Value = 2
$ws.Cells.Item(151, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(151, 3).Value = "Coquimbo"
$ws.Cells.Item(151, 4).Value = 44644
$ws.Cells.Item(151, 5).Value = 4
$ws.Cells.Item(151, 6).Value = 100112031
$ws.Cells.Item(151, 7).Value = "Poroto verde"
$ws.Cells.Item(151, 8).Value = "Magnum"
$ws.Cells.Item(151, 9).Value = "Primera"
$ws.Cells.Item(151, 10).Value = 500
$ws.Cells.Item(151, 11).Value = 20000
$ws.Cells.Item(151, 12).Value = 22000
$ws.Cells.Item(151, 13).Value = 21000
$ws.Cells.Item(151, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(151, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(151, 16).Value = 840
$ws.Cells.Item(151, 17).Value = 25
$ws.Cells.Item(151, 18).Value = "Hortaliza"

# New row 152
$ws.Cells.Item(152, 1).Value = 2
$ws.Cells.Item(152, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(152, 3).Value = "Coquimbo"
$ws.Cells.Item(152, 4).Value = 44644
$ws.Cells.Item(152, 5).Value = 4
$ws.Cells.Item(152, 6).Value = 100112031
$ws.Cells.Item(152, 7).Value = "Poroto verde"
$ws.Cells.Item(152, 8).Value = "Sin especificar"
$ws.Cells.Item(152, 9).Value = "Primera"
$ws.Cells.Item(152, 10).Value = 400
$ws.Cells.Item(152, 11).Value = 24000
$ws.Cells.Item(152, 12).Value = 26000
$ws.Cells.Item(152, 13).Value = 25000
$ws.Cells.Item(152, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(152, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(152, 16).Value = 1000
$ws.Cells.Item(152, 17).Value = 25
$ws.Cells.Item(152, 18).Value = "Hortaliza"
